# Unify the conception of DataNode, DataTable, Entity.
# The sheet that used to represent a generic "Property" table is
# renamed to "DataNode", and the cursor/selection is left where the
# author last left it before saving.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rename the sheet: Property1 -> DataNode
$ws.Name = "DataNode"

# Restore the author's last selection/active cell (C36) on the sheet.
# The sheet view is frozen (panes split at row 8), so this selects the
# cell within the scrollable "bottomLeft" pane, matching the saved file.
$ws.Range("C36").Select() | Out-Null
